# mariculture and workforce update
# Replace the text region-name labels in column A with their numeric rgn_id
# codes, update the jobs figures for the re-mapped / new mariculture regions,
# and append the newly-added region-4 rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 onward: (rgn_id, year, jobs)
$data = @(
    @(1, 2004, 13194),
    @(1, 2005, 13527),
    @(1, 2006, 14606),
    @(1, 2007, 15992),
    @(1, 2008, 17164),
    @(1, 2009, 17627),
    @(1, 2010, 18458),
    @(1, 2011, 18499),
    @(1, 2012, 18927),
    @(1, 2013, 18994),
    @(1, 2014, 20253),
    @(3, 2001, 3391),
    @(3, 2002, 3573),
    @(3, 2003, 3622),
    @(3, 2004, 3726),
    @(3, 2005, 3923),
    @(3, 2006, 3968),
    @(3, 2007, 4025),
    @(3, 2008, 3884),
    @(3, 2009, 3675),
    @(3, 2010, 3724),
    @(3, 2011, 3959),
    @(3, 2012, 3975),
    @(3, 2013, 3910),
    @(3, 2014, 3753),
    @(3, 2015, 3752),
    @(2, 2008, 12222),
    @(2, 2009, 11930),
    @(2, 2010, 13253),
    @(2, 2011, 14089),
    @(2, 2012, 14122),
    @(2, 2013, 14304),
    @(2, 2014, 13782),
    @(8, 2008, 37353.599999999999),
    @(8, 2009, 36721.599999999999),
    @(8, 2010, 35941.599999999999),
    @(8, 2011, 36309.599999999999),
    @(8, 2012, 35992.800000000003),
    @(8, 2013, 35902.400000000001),
    @(8, 2014, 35628),
    @(9, 2008, 9338.4000000000015),
    @(9, 2009, 9180.4000000000015),
    @(9, 2010, 8985.4000000000015),
    @(9, 2011, 9077.4000000000015),
    @(9, 2012, 8998.1999999999971),
    @(9, 2013, 8975.5999999999985),
    @(9, 2014, 8907),
    @(6, 2008, 243681),
    @(6, 2009, 242634),
    @(6, 2010, 243614),
    @(6, 2011, 246256),
    @(6, 2012, 246475),
    @(6, 2013, 249146),
    @(6, 2014, 250526),
    @(5, 2008, 1547),
    @(5, 2009, 1455),
    @(5, 2010, 1411),
    @(5, 2011, 1539),
    @(5, 2012, 1633),
    @(5, 2013, 1594),
    @(5, 2014, 1629),
    @(4, 2010, 1564772.2013915416),
    @(4, 2011, 1561161.9111478101),
    @(4, 2012, 1738006.2274279771),
    @(4, 2013, 1558714.9906560762),
    @(4, 2014, 1545815.9152343646)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$ws.Range("F63").Select()
